$d = $word.ActiveDocument

# 1. Brief paragraph (table cell)
$d.Content.Find.Execute(
    "An email sent to partners in the target country who have sent their documents for review. It will be sent via customer.io",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Một email gửi đến các đối tác ở quốc gia mục tiêu đã gửi tài liệu của họ để xem xét. It will be sent via customer.io",
    2)

# 2. Heading "Thank you for submitting your documents"
$d.Content.Find.Execute(
    "Thank you for submitting your documents",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cảm ơn bạn đã gửi các giấy tờ cần thiết",
    2)

# 3. "Hi " greeting
$d.Content.Find.Execute(
    "Hi ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Xin chào ",
    2)

# 4. "Thank you for providing us with your documents for the upcoming "
$d.Content.Find.Execute(
    "Thank you for providing us with your documents for the upcoming ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cảm ơn bạn đã gửi cho chúng tôi các giấy tờ cần thiết của bạn cho sự kiện ",
    2)

# 5. ". Based on the information you've given us..."
$d.Content.Find.Execute(
    ". Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " sắp tới. Dựa trên thông tin bạn đã cung cấp, chúng tôi sẽ tiến hành sắp xếp chỗ ở và phương tiện đi lại trong quá trình bạn tham gia sự kiện.",
    2)

# 6. "We're currently reviewing your documents..."
$d.Content.Find.Execute(
    "We’re currently reviewing your documents and will reach out to you if we need anything else. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hiện chúng tôi đang kiểm tra giấy tờ của bạn và sẽ liên hệ với bạn nếu chúng tôi cần thêm thông tin. ",
    2)

# --- Paragraph with "If you have any questions, please contact us via <live chat> or <WhatsApp>." ---
# This paragraph has a commentRangeStart immediately before the first run, and the
# run with " or " sits directly between two w:hyperlink elements. A plain Find/Replace
# across the full run text causes the commentRangeStart marker to shift, and causes the
# " or " replacement to inherit the neighbouring hyperlink's character formatting.
# Work around both issues by (a) inserting the new text then removing the old text via a
# second Find (keeps commentRangeStart anchored before the run), and (b) only replacing the
# inner "or" word rather than the whole " or " run (keeps the run's own — absent — rPr).
$pLiveChat = $d.Paragraphs.Item(20)

# 7. "If you have any questions, please contact us via "
$rOld = $pLiveChat.Range.Duplicate
$rOld.Find.Execute("If you have any questions, please contact us via ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rOld.InsertBefore("Nếu bạn cần hỗ trợ, vui lòng liên hệ với chúng tôi qua ")
$rOld2 = $pLiveChat.Range.Duplicate
$rOld2.Find.Execute("If you have any questions, please contact us via ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rOld2.Delete()

# 8. " or " between live chat and WhatsApp hyperlinks -> " hoặc " (replace inner word only)
$pLiveChat.Range.Find.Execute(
    "or", $true, $false, $true, $false, $false, $true, 1, $false,
    "hoặc", 2)

# --- Paragraph with "If you have any questions, please contact your country manager, ..." ---
$pManager = $d.Paragraphs.Item(21)

# 9. "If you have any questions, please contact your country manager, "
$pManager.Range.Find.Execute(
    "If you have any questions, please contact your country manager, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn ",
    2)

# 10. ", at " -> ", qua email "
$pManager.Range.Find.Execute(
    ", at ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", qua email ",
    2)

# 11. " or " between EMAIL ADDRESS and WHATSAPP NO -> " hoặc số "
$pManager.Range.Find.Execute(
    "or", $true, $false, $true, $false, $false, $true, 1, $false,
    "hoặc số", 2)

# 12. "We look forward to seeing you at "
$d.Content.Find.Execute(
    "We look forward to seeing you at ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Chúng tôi rất mong được gặp bạn tại sự kiện ",
    2)
